$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56 ("Yats Fishers") no longer present in the refreshed data -> remove it
$ws.Rows.Item(56).Delete()

$ws.Range("A2").Value = 7
$ws.Range("C2").Value = "10 West Restaurant & Bar"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 4.7
$ws.Range("F2").Value = 692

$ws.Range("A3").Value = 28
$ws.Range("C3").Value = "A2Z Cafe (Inside and patio dining or Carry-out to Curbside)"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 4.7
$ws.Range("F3").Value = 479

$ws.Range("A4").Value = 17
$ws.Range("C4").Value = "Aristocrat Pub & Restaurant"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 4.5
$ws.Range("F4").Value = 1182

$ws.Range("A5").Value = 29
$ws.Range("C5").Value = "Axum Ethiopian Restaurant"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 4.7
$ws.Range("F5").Value = 413

$ws.Range("A6").Value = 57
$ws.Range("C6").Value = "BRU Burger Bar"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 4.6
$ws.Range("F6").Value = 4150

$ws.Range("A7").Value = 11
$ws.Range("C7").Value = "Bluebeard"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 4.7
$ws.Range("F7").Value = 1366

$ws.Range("A8").Value = 24
$ws.Range("C8").Value = "Bonefish Grill"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 4.5
$ws.Range("F8").Value = 1081

$ws.Range("A9").Value = 37
$ws.Range("C9").Value = "Bosphorus Istanbul Cafe"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 4.5
$ws.Range("F9").Value = 1271

$ws.Range("A10").Value = 6
$ws.Range("C10").Value = "Burritos & Beer Restaurant, LLC"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 4.7
$ws.Range("F10").Value = 329

$ws.Range("A11").Value = 31
$ws.Range("C11").Value = "Canal Bistro"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 4.6
$ws.Range("F11").Value = 890

$ws.Range("A12").Value = 33
$ws.Range("C12").Value = "Charleston's Restaurant"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 4.5
$ws.Range("F12").Value = 1058

$ws.Range("A13").Value = 13
$ws.Range("C13").Value = "Chicken Salad Chick"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 4.6
$ws.Range("F13").Value = 30

$ws.Range("A14").Value = 10
$ws.Range("C14").Value = "Chuy's"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 4.4
$ws.Range("F14").Value = 2272

$ws.Range("A15").Value = 5
$ws.Range("C15").Value = "Cooper's Hawk Winery & Restaurant"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 4.6
$ws.Range("F15").Value = 1498

$ws.Range("A16").Value = 54
$ws.Range("C16").Value = "Courses Restaurant"
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = 4.6
$ws.Range("F16").Value = 38

$ws.Range("A17").Value = 36
$ws.Range("C17").Value = "Cracker Barrel Old Country Store"
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 4.4
$ws.Range("F17").Value = 2750

$ws.Range("A18").Value = 50
$ws.Range("C18").Value = "Culver's"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 4.4
$ws.Range("F18").Value = 1563

$ws.Range("A19").Value = 42
$ws.Range("C19").Value = "Fire by the Monon"
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 4.6
$ws.Range("F19").Value = 906

$ws.Range("A20").Value = 46
$ws.Range("C20").Value = "First Watch"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 4.6
$ws.Range("F20").Value = 396

$ws.Range("A21").Value = 56
$ws.Range("C21").Value = "Flatwater"
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 4.6
$ws.Range("F21").Value = 874

$ws.Range("A22").Value = 23
$ws.Range("C22").Value = "Greek Islands"
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 4.6
$ws.Range("F22").Value = 866

$ws.Range("A23").Value = 58
$ws.Range("C23").Value = "His Place Eatery - Chicken & Waffles, Ribs and Soul Food"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 4.5
$ws.Range("F23").Value = 2164

$ws.Range("A24").Value = 22
$ws.Range("C24").Value = "Houlihan's"
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 4.3
$ws.Range("F24").Value = 836

$ws.Range("A25").Value = 27
$ws.Range("C25").Value = "Iron Skillet Restaurant"
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 4.5
$ws.Range("F25").Value = 470

$ws.Range("A26").Value = 34
$ws.Range("C26").Value = "Jamaican Reggae Grill"
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 4.6
$ws.Range("F26").Value = 522

$ws.Range("A27").Value = 20
$ws.Range("C27").Value = "Livery"
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 4.7
$ws.Range("F27").Value = 1490

$ws.Range("A28").Value = 8
$ws.Range("C28").Value = "Maggiano's Little Italy"
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 4.4
$ws.Range("F28").Value = 2257

$ws.Range("A29").Value = 55
$ws.Range("C29").Value = "Mama Carolla's"
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 4.7
$ws.Range("F29").Value = 1639

$ws.Range("A30").Value = 9
$ws.Range("C30").Value = "Meridian Restaurant & Bar"
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 4.5
$ws.Range("F30").Value = 365

$ws.Range("A31").Value = 43
$ws.Range("C31").Value = "Mesh"
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 4.4
$ws.Range("F31").Value = 1150

$ws.Range("A32").Value = 52
$ws.Range("C32").Value = "Mimi Blue Restaurants"
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 4.5
$ws.Range("F32").Value = 666

$ws.Range("A33").Value = 48
$ws.Range("C33").Value = "Nada"
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = 4.4
$ws.Range("F33").Value = 1952

$ws.Range("A34").Value = 30
$ws.Range("C34").Value = "Nesso"
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = 4.7
$ws.Range("F34").Value = 219

$ws.Range("A35").Value = 45
$ws.Range("C35").Value = "O'Charley’s Restaurant & Bar"
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 1336

$ws.Range("A36").Value = 44
$ws.Range("C36").Value = "Ocean Prime"
$ws.Range("D36").Value = 4
$ws.Range("E36").Value = 4.6
$ws.Range("F36").Value = 958

$ws.Range("A37").Value = 15
$ws.Range("C37").Value = "Olive Garden Italian Restaurant"
$ws.Range("D37").Value = 2
$ws.Range("E37").Value = 4.4
$ws.Range("F37").Value = 1395

$ws.Range("A38").Value = 49
$ws.Range("C38").Value = "Pasto Italiano Restaurant & Bar"
$ws.Range("D38").Value = 2
$ws.Range("E38").Value = 4.7
$ws.Range("F38").Value = 195

$ws.Range("A39").Value = 39
$ws.Range("C39").Value = "Perkins Restaurant & Bakery"
$ws.Range("D39").Value = 2
$ws.Range("E39").Value = 4.3
$ws.Range("F39").Value = 981

$ws.Range("A40").Value = 19
$ws.Range("C40").Value = "Ristorante Roma"
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = 4.7
$ws.Range("F40").Value = 159

$ws.Range("A41").Value = 4
$ws.Range("C41").Value = "Rusty Bucket Restaurant and Tavern"
$ws.Range("D41").Value = 2
$ws.Range("E41").Value = 4.4
$ws.Range("F41").Value = 946

$ws.Range("A42").Value = 35
$ws.Range("C42").Value = "Ruth's Chris Steak House"
$ws.Range("D42").Value = 4
$ws.Range("E42").Value = 4.6
$ws.Range("F42").Value = 1709

$ws.Range("A43").Value = 21
$ws.Range("C43").Value = "Sahm's Restaurant"
$ws.Range("D43").Value = 2
$ws.Range("E43").Value = 4.5
$ws.Range("F43").Value = 793

$ws.Range("A44").Value = 3
$ws.Range("C44").Value = "Seasons 52"
$ws.Range("D44").Value = 2
$ws.Range("E44").Value = 4.5
$ws.Range("F44").Value = 1339

$ws.Range("A45").Value = 14
$ws.Range("C45").Value = "Sero's Family Restaurant"
$ws.Range("D45").Value = 2
$ws.Range("E45").Value = 4.5
$ws.Range("F45").Value = 1158

$ws.Range("A46").Value = 51
$ws.Range("C46").Value = "Slapfish"
$ws.Range("D46").Value = 2
$ws.Range("E46").Value = 4.6
$ws.Range("F46").Value = 317

$ws.Range("A47").Value = 47
$ws.Range("C47").Value = "The Bank Restaurant"
$ws.Range("D47").Value = 2
$ws.Range("E47").Value = 4.3
$ws.Range("F47").Value = 455

$ws.Range("A48").Value = 1
$ws.Range("C48").Value = "The Capital Grille"
$ws.Range("D48").Value = 4
$ws.Range("E48").Value = 4.6
$ws.Range("F48").Value = 821

$ws.Range("A49").Value = 41
$ws.Range("C49").Value = "The Cheesecake Factory"
$ws.Range("D49").Value = 2
$ws.Range("E49").Value = 4.2
$ws.Range("F49").Value = 3306

$ws.Range("A50").Value = 59
$ws.Range("C50").Value = "The Italian House on Park"
$ws.Range("D50").Value = 2
$ws.Range("E50").Value = 4.8
$ws.Range("F50").Value = 544

$ws.Range("A51").Value = 18
$ws.Range("C51").Value = "Tinker Street Restaurant"
$ws.Range("D51").Value = 3
$ws.Range("E51").Value = 4.7
$ws.Range("F51").Value = 702

$ws.Range("A52").Value = 38
$ws.Range("C52").Value = "Twin Peaks Restaurant"
$ws.Range("D52").Value = 2
$ws.Range("E52").Value = 4.5
$ws.Range("F52").Value = 3906

$ws.Range("A53").Value = 40
$ws.Range("C53").Value = "Weber Grill Restaurant"
$ws.Range("D53").Value = 2
$ws.Range("E53").Value = 4.2
$ws.Range("F53").Value = 2307

$ws.Range("A54").Value = 16
$ws.Range("C54").Value = "Yard House"
$ws.Range("D54").Value = 2
$ws.Range("E54").Value = 4.4
$ws.Range("F54").Value = 2326

$ws.Range("A55").Value = 12
$ws.Range("C55").Value = "Yats"
$ws.Range("D55").Value = 1
$ws.Range("E55").Value = 4.8
$ws.Range("F55").Value = 1279
